$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New log entry on row 5
$ws.Range("A5").Value = 45617
$ws.Range("A5").NumberFormat = $ws.Range("A4").NumberFormat

$ws.Range("B5").Value = "Movement and dialogue system"

# Selection ends up on B14, as last left by the author
$ws.Range("B14").Select()
